$wb = $excel.ActiveWorkbook

# Rename worksheets (fixed timestamps in sheet names)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777959355574"
$wb.Worksheets.Item(2).Name = "NB_TO-16504777979295547"
$wb.Worksheets.Item(3).Name = "RS_TO-16504777979305549"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504777979925542"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504777980555944"

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504777958945558.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777959195557.csv"
$ws1.Range("B4").Value = "go_stims-16504777959205544.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777959345908.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-16504777966085901.csv"
$ws2.Range("B3").Value = "OB-16504777973905554.csv"
$ws2.Range("B4").Value = "OB-16504777968235853.csv"
$ws2.Range("B5").Value = "ZB-match_2-16504777960575635.csv"
$ws2.Range("B6").Value = "OB-1650477797087555.csv"
$ws2.Range("B7").Value = "TB-1650477797698587.csv"
$ws2.Range("B8").Value = "TB-16504777974625874.csv"
$ws2.Range("B9").Value = "ZB-match_5-16504777967395873.csv"
$ws2.Range("B10").Value = "TB-16504777979105525.csv"

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777979455562.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777979325566.csv"
$ws4.Range("B4").Value = "MM_stims-1650477797976587.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477797946557.csv"
$ws4.Range("B6").Value = "MM_stims-16504777979915528.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477797976587.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504777979955554.csv"
$ws5.Range("B3").Value = "vSAT_stims-1650477798039553.csv"
$ws5.Range("B4").Value = "SAT_stims-16504777980075576.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477798023555.csv"
